$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text interpretation so that
# numeric-looking strings (e.g. "601.58") are not coerced into numbers,
# matching the original inline-string cell typing used throughout the sheet.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '68.354.73'
$ws.Range("E2").Value = '  -1.08%  '
Set-TextValue $ws.Range("D3") '3.891.35'
$ws.Range("E3").Value = '  +2.19%  '
$ws.Range("E4").Value = '  +0.13%  '
Set-TextValue $ws.Range("D5") '601.58'
$ws.Range("E5").Value = '  +0.00%  '
Set-TextValue $ws.Range("D6") '166.85'
$ws.Range("E6").Value = '  +1.79%  '
Set-TextValue $ws.Range("D7") '3.887.49'
$ws.Range("E7").Value = '  +2.15%  '
$ws.Range("E8").Value = '  +0.00%  '
Set-TextValue $ws.Range("D9") '0.527'
$ws.Range("E9").Value = '  -1.54%  '
Set-TextValue $ws.Range("D10") '0.167'
$ws.Range("E10").Value = '  -1.98%  '
Set-TextValue $ws.Range("D11") '6.43'
$ws.Range("E11").Value = '  +2.03%  '
Set-TextValue $ws.Range("D12") '0.460'
$ws.Range("E12").Value = '  -0.46%  '
Set-TextValue $ws.Range("D13") '0.0000254'
$ws.Range("E13").Value = '  +3.17%  '
Set-TextValue $ws.Range("D14") '37.34'
$ws.Range("E14").Value = '  +0.23%  '
Set-TextValue $ws.Range("D15") '4.554.42'
$ws.Range("E15").Value = '  +2.47%  '
Set-TextValue $ws.Range("D16") '3.908.13'
$ws.Range("E16").Value = '  +2.67%  '
Set-TextValue $ws.Range("D17") '68.512.68'
$ws.Range("E17").Value = '  -1.02%  '
Set-TextValue $ws.Range("D18") '7.45'
$ws.Range("E18").Value = '  +0.12%  '
Set-TextValue $ws.Range("D19") '17.26'
$ws.Range("E19").Value = '  -0.49%  '
Set-TextValue $ws.Range("D20") '0.111'
$ws.Range("E20").Value = '  -2.14%  '
Set-TextValue $ws.Range("D21") '11.03'
$ws.Range("E21").Value = '  -3.38%  '
Set-TextValue $ws.Range("D22") '488.39'
$ws.Range("E22").Value = '  +0.05%  '
Set-TextValue $ws.Range("D23") '0.725'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  +4.09%  '
Set-TextValue $ws.Range("D25") '84.58'
$ws.Range("E25").Value = '  -0.31%  '
Set-TextValue $ws.Range("D26") '2.23'
$ws.Range("E26").Value = '  -1.53%  '
Set-TextValue $ws.Range("D27") '12.00'
$ws.Range("E27").Value = '  -1.92%  '
Set-TextValue $ws.Range("D28") '10.13'
$ws.Range("E28").Value = '  +0.96%  '
$ws.Range("E29").Value = '  -0.02%  '
Set-TextValue $ws.Range("D30") '2.93'
$ws.Range("E30").Value = '  -1.36%  '
Set-TextValue $ws.Range("D31") '4.043.80'
$ws.Range("E31").Value = '  +2.17%  '
Set-TextValue $ws.Range("D32") '2.37'
$ws.Range("E32").Value = '  -1.41%  '
Set-TextValue $ws.Range("D33") '7.71'
$ws.Range("E33").Value = '  -3.90%  '
Set-TextValue $ws.Range("D34") '31.74'
$ws.Range("E34").Value = '  -0.22%  '
Set-TextValue $ws.Range("D35") '3.848.84'
$ws.Range("E35").Value = '  +2.57%  '
Set-TextValue $ws.Range("D36") '0.106'
$ws.Range("E36").Value = '  -0.68%  '
$ws.Range("E37").Value = '  +1.20%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D38") '0.139'
$ws.Range("E38").Value = '  -1.58%  '
$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D39") '5.92'
$ws.Range("E39").Value = '  +0.18%  '
Set-TextValue $ws.Range("D40") '3.16'
$ws.Range("E40").Value = '  +3.97%  '
$ws.Range("E41").Value = '  +0.12%  '
Set-TextValue $ws.Range("D42") '0.316'
$ws.Range("E42").Value = '  -1.26%  '
Set-TextValue $ws.Range("D43") '432.06'
$ws.Range("E43").Value = '  +1.44%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D44") '48.21'
$ws.Range("E44").Value = '  -0.82%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D45") '1.98'
$ws.Range("E45").Value = '  -0.71%  '
Set-TextValue $ws.Range("D46") '8.50'
$ws.Range("E46").Value = '  +1.37%  '
Set-TextValue $ws.Range("D48") '142.52'
$ws.Range("E48").Value = '  +0.83%  '
$ws.Range("E49").Value = '  +17.67%  '
Set-TextValue $ws.Range("D50") '2.796.92'
$ws.Range("E50").Value = '  -1.36%  '
Set-TextValue $ws.Range("D51") '39.28'
$ws.Range("E51").Value = '  -0.76%  '
